# Add a new data row to the "matriz" sheet describing the
# sia-estacio system's professor/aluno profile conflict.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("matriz")

$ws.Range("A3").Value = "sia-estacio"
$ws.Range("B3").Value = "professor"
$ws.Range("C3").Value = "sia-estacio"
$ws.Range("D3").Value = "aluno"
